# Script is working till company creating
#
# Adds a new "GitHubSync" worksheet right after "AddCustomerTest" (and
# before "OpenAccountTest"), populates it with the username/companyname
# header row plus one data row, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("AddCustomerTest")
$gitHubSync = $wb.Worksheets.Add($null, $afterSheet)
$gitHubSync.Name = "GitHubSync"

$gitHubSync.Range("A1").Value = "username"
$gitHubSync.Range("B1").Value = "companyname"
$gitHubSync.Range("B2").Value = "vase123"
$gitHubSync.Range("A2").Value = "bhautik-vase45"

[void]$gitHubSync.Activate()
[void]$gitHubSync.Range("A3").Select()
